# Append the new daily COVID-19 data row (14 June 2020 update -> date 13 June 2020,
# serial 43995) to the "Tabela1" table on the "Covid-19 podatki" sheet.
#
# The table currently spans A1:J94 (header in row 1, data in rows 2-94).
# We grow the table by one row (ListRows.Add) so the table reference,
# the AutoFilter reference and the sheet dimension all extend to J95
# automatically, then copy the formatting of the previous data row (93)
# onto the freshly added row so the banded/bordered look is preserved,
# and finally fill in the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Copy the format of the last-but-one data row so the new row picks up the
# same font / number format / border as its neighbours, then grow the table.
$ws.Range("A93:J93").Copy() | Out-Null
$newRow = $lo.ListRows.Add()
$ws.Range("A95:J95").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New data values for 13/6/2020 (date serial 43995).
$ws.Range("A95").Value = 43995
$ws.Range("B95").Value = 87386
$ws.Range("C95").Value = 291
$ws.Range("D95").Value = 1495
$ws.Range("E95").Value = 3
$ws.Range("F95").Value = 6
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 109
$ws.Range("J95").Value = 0

# Match the saved selection (whole new row selected, active cell A95).
$ws.Range("A95:J95").Select() | Out-Null

Write-Host "Added row 95 to Tabela1; table now spans" $lo.Range.Address()
